# Apply the price / volume(1h) refresh captured by the commit
# "Updated cryptos list on Mon Aug 21 15:11:18 UTC 2023 with GitHub Actions".
# Column D ("Price") and E ("Volume(1h)") are plain text cells in the source
# workbook (t="inlineStr"), and rows 48/49 (EnergySwap/Cronos) swap order too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.192.95"
$ws.Cells.Item(2, 5).Value = "  -0.54%  "
$ws.Cells.Item(3, 4).Value = "1.678.85"
$ws.Cells.Item(3, 5).Value = "  -0.44%  "
$ws.Cells.Item(4, 4).Value = "'1.005"
$ws.Cells.Item(4, 5).Value = "  -0.85%  "
$ws.Cells.Item(5, 4).Value = "'210.08"
$ws.Cells.Item(5, 5).Value = "  -3.82%  "
$ws.Cells.Item(6, 4).Value = "'0.5283"
$ws.Cells.Item(6, 5).Value = "  -2.85%  "
$ws.Cells.Item(7, 4).Value = "'1.005"
$ws.Cells.Item(7, 5).Value = "  -0.79%  "
$ws.Cells.Item(8, 4).Value = "'0.2675"
$ws.Cells.Item(8, 5).Value = "  -1.48%  "
$ws.Cells.Item(9, 4).Value = "'0.06289"
$ws.Cells.Item(9, 5).Value = "  -2.45%  "
$ws.Cells.Item(10, 4).Value = "'21.33"
$ws.Cells.Item(10, 5).Value = "  -3.02%  "
$ws.Cells.Item(11, 4).Value = "'0.07534"
$ws.Cells.Item(11, 5).Value = "  -2.07%  "
$ws.Cells.Item(12, 4).Value = "1.676.44"
$ws.Cells.Item(12, 5).Value = "  -0.80%  "
$ws.Cells.Item(13, 4).Value = "'4.466"
$ws.Cells.Item(13, 5).Value = "  -1.30%  "
$ws.Cells.Item(14, 4).Value = "'0.5657"
$ws.Cells.Item(14, 5).Value = "  -2.45%  "
$ws.Cells.Item(15, 4).Value = "'0.000008121"
$ws.Cells.Item(15, 5).Value = "  -2.54%  "
$ws.Cells.Item(16, 4).Value = "'66.46"
$ws.Cells.Item(16, 5).Value = "  +2.19%  "
$ws.Cells.Item(17, 4).Value = "26.243.75"
$ws.Cells.Item(17, 5).Value = "  -0.71%  "
$ws.Cells.Item(18, 4).Value = "'1.004"
$ws.Cells.Item(18, 5).Value = "  -0.77%  "
$ws.Cells.Item(19, 4).Value = "'4.849"
$ws.Cells.Item(19, 5).Value = "  -1.98%  "
$ws.Cells.Item(20, 4).Value = "'10.51"
$ws.Cells.Item(20, 5).Value = "  -4.06%  "
$ws.Cells.Item(21, 4).Value = "'188.22"
$ws.Cells.Item(21, 5).Value = "  -0.97%  "
$ws.Cells.Item(22, 4).Value = "'6.208"
$ws.Cells.Item(22, 5).Value = "  -0.09%  "
$ws.Cells.Item(23, 5).Value = "  -0.75%  "
$ws.Cells.Item(24, 4).Value = "'146.89"
$ws.Cells.Item(24, 5).Value = "  -2.21%  "
$ws.Cells.Item(25, 4).Value = "'0.1257"
$ws.Cells.Item(25, 5).Value = "  -3.23%  "
$ws.Cells.Item(26, 4).Value = "'7.627"
$ws.Cells.Item(26, 5).Value = "  -3.01%  "
$ws.Cells.Item(27, 4).Value = "'15.90"
$ws.Cells.Item(27, 5).Value = "  +1.18%  "
$ws.Cells.Item(28, 4).Value = "'0.06411"
$ws.Cells.Item(28, 5).Value = "  +1.22%  "
$ws.Cells.Item(29, 4).Value = "'1.345"
$ws.Cells.Item(29, 5).Value = "  -4.70%  "
$ws.Cells.Item(30, 4).Value = "'1.280"
$ws.Cells.Item(30, 5).Value = "  -3.54%  "
$ws.Cells.Item(31, 4).Value = "'3.537"
$ws.Cells.Item(31, 5).Value = "  -0.77%  "
$ws.Cells.Item(32, 4).Value = "'3.483"
$ws.Cells.Item(32, 5).Value = "  -2.48%  "
$ws.Cells.Item(33, 4).Value = "'1.652"
$ws.Cells.Item(33, 5).Value = "  -1.06%  "
$ws.Cells.Item(34, 4).Value = "'1.013"
$ws.Cells.Item(34, 5).Value = "  -2.40%  "
$ws.Cells.Item(35, 4).Value = "'0.6080"
$ws.Cells.Item(35, 5).Value = "  -2.00%  "
$ws.Cells.Item(36, 4).Value = "'2.413"
$ws.Cells.Item(36, 5).Value = "  -0.24%  "
$ws.Cells.Item(37, 4).Value = "'2.718"
$ws.Cells.Item(37, 5).Value = "  -0.22%  "
$ws.Cells.Item(38, 4).Value = "'6.164"
$ws.Cells.Item(38, 5).Value = "  -0.86%  "
$ws.Cells.Item(39, 4).Value = "'0.01613"
$ws.Cells.Item(39, 5).Value = "  -1.31%  "
$ws.Cells.Item(40, 4).Value = "1.099.86"
$ws.Cells.Item(40, 5).Value = "  -1.34%  "
$ws.Cells.Item(41, 4).Value = "'0.8693"
$ws.Cells.Item(41, 5).Value = "  -1.21%  "
$ws.Cells.Item(42, 4).Value = "'1.008"
$ws.Cells.Item(42, 5).Value = "  -0.93%  "
$ws.Cells.Item(43, 4).Value = "'100.01"
$ws.Cells.Item(43, 5).Value = "  -1.05%  "
$ws.Cells.Item(44, 4).Value = "1.833.00"
$ws.Cells.Item(44, 5).Value = "  -0.44%  "
$ws.Cells.Item(45, 4).Value = "'0.00000000110"
$ws.Cells.Item(45, 5).Value = "  -0.10%  "
$ws.Cells.Item(46, 4).Value = "'56.86"
$ws.Cells.Item(46, 5).Value = "  -0.66%  "
$ws.Cells.Item(47, 5).Value = "  -0.95%  "
$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48, 4).Value = "'0.05257"
$ws.Cells.Item(48, 5).Value = "  -0.37%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).Value = "'7.999"
$ws.Cells.Item(49, 5).Value = "  -2.57%  "
$ws.Cells.Item(50, 4).Value = "'0.4269"
$ws.Cells.Item(50, 5).Value = "  -0.80%  "
$ws.Cells.Item(51, 4).Value = "'5.958"
$ws.Cells.Item(51, 5).Value = "  -1.47%  "
